# course bug(elective, MC major) fix
# Update the "DB Updated Date" value in cell B1 to reflect the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 holds the "DB Updated Date" as text (e.g. "2021.02.18"). Force the
# number format to Text before writing so Excel doesn't auto-convert the
# dotted date string into a date serial, then restore the default format
# so the cell keeps its original (unstyled) appearance.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "2021.02.25"
$ws.Range("B1").ClearFormats()
